# Insert a new weekly price record as row 13 in the "Albahaca" sheet.
# This pushes the existing rows 13..102 down to 14..103 (dimension grows
# from A1:R102 to A1:R103), and populates the newly inserted row 13 with
# the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 13, shifting rows 13-102 down
# to 14-103 (this also extends the sheet dimension automatically).
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new record's values.
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44761
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112052
$ws.Range("G13").Value = "Albahaca"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 1400
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4500
$ws.Range("M13").Value = 4250
$ws.Range("N13").Value = "$/paquete"
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 4250
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
